# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (exhibit listing) ---
$wsExhibit.Range("G4").Value = "不可售"
$wsExhibit.Range("F7").Value = 1409
$wsExhibit.Range("F8").Value = 216
$wsExhibit.Range("F9").Value = 76
$wsExhibit.Range("F10").Value = 115
$wsExhibit.Range("F11").Value = 6031
$wsExhibit.Range("F12").Value = 63
$wsExhibit.Range("F15").Value = 4822
$wsExhibit.Range("F17").Value = 172
$wsExhibit.Range("F18").Value = 1158
$wsExhibit.Range("F19").Value = 47
$wsExhibit.Range("F23").Value = 275
$wsExhibit.Range("F25").Value = 3296
$wsExhibit.Range("F26").Value = 136

# --- Sheet "全部类型" (all types, combined listing) ---
$wsAll.Range("G4").Value = "不可售"
$wsAll.Range("F8").Value = 1409
$wsAll.Range("F9").Value = 216
$wsAll.Range("F10").Value = 76
$wsAll.Range("F11").Value = 115
$wsAll.Range("F12").Value = 6031
$wsAll.Range("F13").Value = 63
$wsAll.Range("F16").Value = 4822
$wsAll.Range("F18").Value = 172
$wsAll.Range("F19").Value = 1158
$wsAll.Range("F20").Value = 47
$wsAll.Range("F24").Value = 275
$wsAll.Range("F26").Value = 3296
$wsAll.Range("F28").Value = 136
